$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Videnov"
$ws.Range("B8").Value = "Sofia, Tsarigradsko, 15"
$ws.Range("A9").Value = "Videnov"
$ws.Range("B9").Value = "Sofia, Tsarigradsko, 15"
